$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price + Volume(1h)) per diff
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.303.94"
$ws.Range("E2").Value = "  +3.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.634.71"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.37"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.08"
$ws.Range("E6").Value = "  +4.51%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").Value = "  +1.68%  "
$ws.Range("E9").Value = "  +8.23%  "
$ws.Range("E10").Value = "  +4.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.76"
$ws.Range("E11").Value = "  +2.32%  "
$ws.Range("E12").Value = "  +1.81%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.92"
$ws.Range("E13").Value = "  +5.85%  "
$ws.Range("E14").Value = "  +19.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.107.05"
$ws.Range("E15").Value = "  +2.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.157.42"
$ws.Range("E16").Value = "  +3.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.641.93"
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("E18").Value = "  +2.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.82"
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "354.87"
$ws.Range("E20").Value = "  +3.06%  "
$ws.Range("E21").Value = "  +5.68%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.95"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.70"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.48"
$ws.Range("E25").Value = "  +3.05%  "
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0949"
$ws.Range("E29").Value = "  +12.89%  "
$ws.Range("E30").Value = "  -0.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.12"
$ws.Range("E31").Value = "  +4.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "517.08"
$ws.Range("E32").Value = "  -6.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.78"
$ws.Range("E33").Value = "  +3.00%  "
$ws.Range("E34").Value = "  +9.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.36"
$ws.Range("E35").Value = "  +4.87%  "
$ws.Range("E36").Value = "  +4.33%  "
$ws.Range("E37").Value = "  +6.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "165.25"
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("E39").Value = "  +4.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("E41").Value = "  +0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.16"
$ws.Range("E42").Value = "  +6.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "164.91"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("E44").Value = "  +4.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0610"
$ws.Range("E45").Value = "  +5.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.09"
$ws.Range("E46").Value = "  +3.91%  "
$ws.Range("E47").Value = "  +10.00%  "
$ws.Range("E48").Value = "  +3.28%  "
$ws.Range("E49").Value = "  +2.05%  "
$ws.Range("E50").Value = "  +2.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.51"
$ws.Range("E51").Value = "  +3.04%  "
